$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fullname in row 2 (C2) to append a digit, simulating invalid data
# used to validate that fullname should not contain digits
$ws.Range("C2").Value = "Phạm Thanh Hà0"

# Update the active selection to reflect the edited cell
$ws.Range("C2").Select()
